# Automatic update of files.
# Applies the row-level data corrections/swaps to rows 52-67 in the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A52").Value = 111901585
$ws.Range("Q52").Value = 478339
$ws.Range("R52").Value = 7035076
$ws.Range("AC52").Value = "ringhack äldre"
$ws.Range("A53").Value = 111901584
$ws.Range("Q53").Value = 478211
$ws.Range("R53").Value = 7035067
$ws.Range("AC53").Value = "ringhack"
$ws.Range("B54").Value = 85197
$ws.Range("B55").Value = 89047
$ws.Range("A56").Value = 111901519
$ws.Range("B56").Value = 86371
$ws.Range("E56").Value = 4412
$ws.Range("F56").Value = "Äggvaxskivling"
$ws.Range("G56").Value = "Hygrophorus karstenii"
$ws.Range("H56").Value = "Sacc. & Cub."
$ws.Range("Q56").Value = 477765
$ws.Range("R56").Value = 7033404
$ws.Range("AC56").Value = ""
$ws.Range("A57").Value = 111901546
$ws.Range("B57").Value = 56430
$ws.Range("E57").Value = 100109
$ws.Range("F57").Value = "Tretåig hackspett"
$ws.Range("G57").Value = "Picoides tridactylus"
$ws.Range("H57").Value = "(Linnaeus, 1758)"
$ws.Range("I57").Value = ""
$ws.Range("N57").Value = ""
$ws.Range("Q57").Value = 477668
$ws.Range("R57").Value = 7033374
$ws.Range("AC57").Value = "ringhack äldre"
$ws.Range("A58").Value = 111901547
$ws.Range("B58").Value = 56430
$ws.Range("E58").Value = 100109
$ws.Range("F58").Value = "Tretåig hackspett"
$ws.Range("G58").Value = "Picoides tridactylus"
$ws.Range("H58").Value = "(Linnaeus, 1758)"
$ws.Range("Q58").Value = 477524
$ws.Range("R58").Value = 7033330
$ws.Range("AC58").Value = "ringhack"
$ws.Range("A59").Value = 111901551
$ws.Range("Q59").Value = 477433
$ws.Range("R59").Value = 7033429
$ws.Range("AC59").Value = "ringhack"
$ws.Range("A60").Value = 111901518
$ws.Range("B60").Value = 86371
$ws.Range("E60").Value = 4412
$ws.Range("F60").Value = "Äggvaxskivling"
$ws.Range("G60").Value = "Hygrophorus karstenii"
$ws.Range("H60").Value = "Sacc. & Cub."
$ws.Range("Q60").Value = 477674
$ws.Range("R60").Value = 7033500
$ws.Range("AC60").Value = ""
$ws.Range("A61").Value = 111901587
$ws.Range("B61").Value = 56575
$ws.Range("E61").Value = 103021
$ws.Range("F61").Value = "Talltita"
$ws.Range("G61").Value = "Poecile montanus"
$ws.Range("H61").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("I61").Value = "'2"
$ws.Range("N61").Value = "observerad"
$ws.Range("Q61").Value = 477611
$ws.Range("R61").Value = 7033311
$ws.Range("AC61").Value = ""
$ws.Range("A62").Value = 111901548
$ws.Range("B62").Value = 56430
$ws.Range("E62").Value = 100109
$ws.Range("F62").Value = "Tretåig hackspett"
$ws.Range("G62").Value = "Picoides tridactylus"
$ws.Range("H62").Value = "(Linnaeus, 1758)"
$ws.Range("Q62").Value = 477476
$ws.Range("R62").Value = 7033385
$ws.Range("AC62").Value = "ringhack äldre"
$ws.Range("A63").Value = 111901618
$ws.Range("B63").Value = 85197
$ws.Range("E63").Value = 249278
$ws.Range("F63").Value = "Barrviolspindling"
$ws.Range("G63").Value = "Cortinarius harcynicus"
$ws.Range("H63").Value = "(Pers.) M.M.Moser"
$ws.Range("Q63").Value = 477471
$ws.Range("R63").Value = 7033412
$ws.Range("AC63").Value = ""
$ws.Range("A64").Value = 111901549
$ws.Range("Q64").Value = 477464
$ws.Range("R64").Value = 7033364
$ws.Range("AC64").Value = "ringhack färska"
$ws.Range("A65").Value = 111901544
$ws.Range("Q65").Value = 477639
$ws.Range("R65").Value = 7033515
$ws.Range("A66").Value = 111901550
$ws.Range("B66").Value = 56430
$ws.Range("E66").Value = 100109
$ws.Range("F66").Value = "Tretåig hackspett"
$ws.Range("G66").Value = "Picoides tridactylus"
$ws.Range("H66").Value = "(Linnaeus, 1758)"
$ws.Range("Q66").Value = 477473
$ws.Range("R66").Value = 7033404
$ws.Range("AC66").Value = "ringhack äldre"
$ws.Range("A67").Value = 111901545
$ws.Range("Q67").Value = 477667
$ws.Range("R67").Value = 7033500
